# Update countries & provincias Spain
#
# The underlying dataset is sorted descending by "Casos totales" (column B).
# A handful of countries' case counts were refreshed, which changed their
# rank and therefore swapped their row with a neighbouring country. For each
# affected row we write out the country name (column A) together with the
# refreshed statistics (columns B-H) so the sheet ends up in the same state
# as after the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# $Row, $Country, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes
function Set-Row($Row, $Country, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# Australia: refreshed stats only, no reordering.
Set-Row 50 "Australia" 6822 21 5849 878 29 0 95

# Bolivia overtakes Guinea / Republica de Macedonia / Nueva Zelanda.
Set-Row 81 "Bolivia" 1594 124 166 1352 3 5 76
Set-Row 82 "Guinea" 1586 0 405 1174 0 0 7
Set-Row 83 "Republica de Macedonia" 1511 0 945 482 21 0 84
Set-Row 84 "Nueva Zelanda" 1487 0 1276 191 0 0 20

# Honduras overtakes Hong Kong / Tunez.
Set-Row 91 "Honduras" 1055 45 118 855 10 6 82
Set-Row 92 "Hong Kong" 1040 0 879 157 3 0 4
Set-Row 93 "Tunez" 1013 0 328 643 25 0 42

# Barbados overtakes Liechtenstein.
Set-Row 162 "Barbados" 82 1 44 31 4 0 7
Set-Row 163 "Liechtenstein" 82 0 55 26 0 0 1

# Belice overtakes Santa Lucia.
Set-Row 188 "Belice" 18 0 13 3 1 0 2
Set-Row 189 "Santa Lucia" 18 0 15 3 0 0 0

# San Vicente y las Granadinas overtakes Namibia (identical totals, order swaps).
Set-Row 194 "San Vicente y las Granadinas" 16 0 8 8 0 0 0
Set-Row 195 "Namibia" 16 0 8 8 0 0 0

# San Cristobal y Nieves overtakes Burundi.
Set-Row 198 "San Cristobal y Nieves" 15 0 8 7 0 0 0
Set-Row 199 "Burundi" 15 0 7 7 0 0 1
